# Update the LSIP data table (221004DataTable) rows for:
#  - Key Stage 4 (KS4) destinations
#  - Key Stage 5 (KS5) destinations
#  - Job adverts by profession
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11: Key Stage 4 (KS4) destinations
$ws.Range("A11").Value = "Key Stage 4 (KS4) destinations  - provisional"
$ws.Range("B11").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/75e2be32-3c51-4790-2c28-08dab0fa305d'>National Pupil Database</a>"
$ws.Range("C11").Value = "Aug 2020 -  Jul 2021 (19/20 learners) (20/10/22)"
$ws.Range("D11").Value = "February 2023 - revision"
$ws.Range("D11").Style = "Normal"

# Row 12: Key Stage 5 (KS5) destinations
$ws.Range("A12").Value = "Key Stage 5 (KS5) destinations - provisional"
$ws.Range("B12").Value = "<a href = 'https://explore-education-statistics.service.gov.uk/data-tables/permalink/62b04091-a13b-40e9-52d9-08dab0fd4449'>National Pupil Database</a>"
$ws.Range("C12").Value = "Aug 2020 -  Jul 2021 (19/20 learners) (20/10/22)"
$ws.Range("D12").Value = "February 2023 - revision"
$ws.Range("D12").Style = "Normal"

# Row 13: Job adverts by profession
$ws.Range("B13").Value = "<a href='https://www.ons.gov.uk/employmentandlabourmarket/peopleinwork/employmentandemployeetypes/datasets/labourdemandvolumesbyprofessionandlocalauthorityuk'>ONS Textkernel</a>"
$ws.Range("C13").Value = "Oct 2022 (21/12/22)"

# Leave the selection where the author last left it before saving
$ws.Range("A5").Select()
